# Fixed Asset Depreciation Code change
#
# Replaces the first four test-case rows (rows 2-5) on every sheet in the
# workbook: the old "JobCreation / CreateExpenses|CreateBudget / Timesheet /
# EmployeeCreation" rows are swapped out for the
# "QueryAndValidateExistingJob / QueryAndValidateExistingEmployee /
# PrintJobBudgetMPL / VerfiyTimesheet" rows (rows 6-8 of GlobalTestPack and
# the tail rows of CHN_SysTest_FullCycle keep their original content).

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "CHN_SysTest_FullCycle",
    "IND_SysTest_FullCycle",
    "IND_Regression_FullCycle",
    "CHN_Regression_FullCycle",
    "Smoke",
    "GlobalTestPack"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A2").Value = "QueryAndValidateExistingJob"
    $ws.Range("B2").Value = "QueryExistingJob"
    $ws.Range("C2").Value = "QueryTheExistingJob"
    $ws.Range("D2").Value = "No"

    $ws.Range("A3").Value = "QueryAndValidateExistingEmployee"
    $ws.Range("B3").Value = "QueryExistingEmployee"
    $ws.Range("C3").Value = "QueryExistingEmployee"
    $ws.Range("D3").Value = "No"

    $ws.Range("A4").Value = "PrintJobBudgetMPL"
    $ws.Range("B4").Value = "JobBudgetMPL"
    $ws.Range("C4").Value = "JobBudgetMPL"
    $ws.Range("D4").Value = "Yes"

    $ws.Range("A5").Value = "VerfiyTimesheet"
    $ws.Range("B5").Value = "verifytimesheet"
    $ws.Range("C5").Value = "verifytimesheet"
    $ws.Range("D5").Value = "No"
}

# Update the selection / active cell on each sheet to match the post-edit
# state, leaving CHN_SysTest_FullCycle selected/active last.
$wb.Worksheets.Item("IND_SysTest_FullCycle").Range("A2:D5").Select()
$wb.Worksheets.Item("IND_Regression_FullCycle").Range("A7").Select()
$wb.Worksheets.Item("CHN_Regression_FullCycle").Range("A2:D5").Select()
$wb.Worksheets.Item("Smoke").Range("A2:D5").Select()
$wb.Worksheets.Item("GlobalTestPack").Range("A7").Select()

$ws1 = $wb.Worksheets.Item("CHN_SysTest_FullCycle")
$ws1.Activate()
$ws1.Range("A8").Select()
